$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 126 to 130 (Heap) completed: mark column C "Done" for rows 336-353
# as "yes". The sheet has been using a rotating set of colored-fill styles
# (s=11..17) for every other completed row instead of one fixed style, so
# reproduce that by copying the cell format from a donor row that already
# carries the wanted style, then stamping the new value in afterwards.
$doneRowStyleDonor = [ordered]@{
    336 = 6
    337 = 6
    338 = 13
    339 = 22
    340 = 22
    341 = 20
    342 = 13
    343 = 13
    344 = 13
    345 = 20
    346 = 16
    347 = 16
    348 = 13
    349 = 6
    350 = 13
    351 = 6
    352 = 13
    353 = 20
}

foreach ($row in $doneRowStyleDonor.Keys) {
    $donorRow = $doneRowStyleDonor[$row]
    $ws.Range("C$donorRow").Copy()
    $ws.Range("C$row").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("C$row").Value = "yes"
}

$excel.CutCopyMode = 0

# Leave the sheet scrolled/selected where the user ended up after
# finishing this range.
$null = $ws.Range("C350").Select()
